$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 732
$ws1.Range("F5").Value = 2790
$ws1.Range("F7").Value = 3734
$ws1.Range("F9").Value = 945
$ws1.Range("F10").Value = 15

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 732
$ws4.Range("F6").Value = 2790
$ws4.Range("F8").Value = 3734
$ws4.Range("F10").Value = 945
$ws4.Range("F11").Value = 15
